# Atualização de bases das ligas, do dia: 03-04-2024 às 22:09
# Two pairs of match rows were re-sorted (their id/Div/Date stayed put but all
# other match data swapped), and a handful of odds were refreshed on several
# still-upcoming fixtures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($rowA, $rowB) {
    # Column B (match id) swaps independently of the F:AC block because C:E
    # (Div / Div Original Name / Date) must stay untouched.
    $bA = $ws.Range("B$rowA").Value2
    $bB = $ws.Range("B$rowB").Value2
    $ws.Range("B$rowA").Value2 = $bB
    $ws.Range("B$rowB").Value2 = $bA

    $rangeA = $ws.Range("F${rowA}:AC${rowA}")
    $rangeB = $ws.Range("F${rowB}:AC${rowB}")
    $valsA = $rangeA.Value2
    $valsB = $rangeB.Value2
    $rangeA.Value2 = $valsB
    $rangeB.Value2 = $valsA
}

Swap-Rows 161 162
Swap-Rows 166 167

# Refreshed odds on upcoming fixtures
$ws.Range("R168").Value2 = 1.92
$ws.Range("S168").Value2 = 1.98
$ws.Range("U168").Value2 = 2
$ws.Range("V168").Value2 = 1.85

$ws.Range("R169").Value2 = 1.9
$ws.Range("S169").Value2 = 2

$ws.Range("U170").Value2 = 1.875
$ws.Range("V170").Value2 = 1.975

$ws.Range("R172").Value2 = 2.05
$ws.Range("S172").Value2 = 1.85

$ws.Range("O173").Value2 = 3.6
$ws.Range("R173").Value2 = 1.83
$ws.Range("S173").Value2 = 2.07
